# Replace the "Our goal is to create ... afterwards." paragraph (project goal
# summary) with its new Bulgarian text, split across runs so that the
# plain-language filler ("platforma ... vremeto" and final ".") that
# originally came from separate runs keeps separate run boundaries while
# the rest carries an explicit bg-BG language tag.
$d = $word.ActiveDocument

$find = $d.Content
$found = $find.Find.Execute("Our goal is to create", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw 'Could not locate the target paragraph (Our goal is to create...)'
}

$para = $find.Paragraphs(1)
$target = $para.Range

$target.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="426"/><w:rPr><w:lang w:val="bg-BG"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">Целта на проекта ми е да създам </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>платформа за управление на проекти и проследяване на времето</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>за работа</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve"> Потребителят ще може да се регистрира с имейла и паролата си и да види задачите, които са му зададени или да създаде нови. Също така ще може да редактира или изтрива вече създадени задачи. </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">Трябва да </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">има създаден админски акаунт с различни </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t xml:space="preserve">от на другите акаунти </w:t></w:r><w:r><w:rPr><w:lang w:val="bg-BG"/></w:rPr><w:t>привилегии като например да може да създава, изтрива и редактира потребители.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
